# Swap the order of slides 5 and 6.
#
# Currently:
#   position 5 -> "Lab #2: Syscall #4 (print string)" slide
#   position 6 -> "Lab #2: Syscalls" slide
#
# After this script:
#   position 5 -> "Lab #2: Syscalls" slide
#   position 6 -> "Lab #2: Syscall #4 (print string)" slide

$p = $ppt.ActivePresentation

# Move the slide currently at position 6 to position 5; PowerPoint shifts
# the slide that was at 5 down to 6, giving the swapped order.
$s = $p.Slides.Item(6)
$s.MoveTo(5)
